$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Select E8 to mirror the active selection recorded in the saved file
$ws.Activate()
$ws.Range("E8").Select()

# Replace the "Good Morning" text with "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"
